$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from H1 to I1/J1 (header style), then set new header values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
